$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12608
$ws1.Range("F5").Value = 34
$ws1.Range("F7").Value = 11
$ws1.Range("F8").Value = 12508
$ws1.Range("F9").Value = 250
$ws1.Range("F10").Value = 4916
$ws1.Range("F11").Value = 4839

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12608
$ws4.Range("F6").Value = 34
$ws4.Range("F8").Value = 11
$ws4.Range("F9").Value = 12508
$ws4.Range("F10").Value = 250
$ws4.Range("F11").Value = 4916
$ws4.Range("F12").Value = 4839
